# Rows 8-23 of the sheet (data rows) get re-shuffled: every row's entire
# content moves to a different row position (row 12 stays put). Read the
# whole A8:AY23 block into memory first so the permutation can be applied
# without any row clobbering another before it's been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A8:AY23")
$data = $srcRange.Value2

# Map: old row number -> new row number (identity for row 12)
$rowMap = @{
    8  = 10
    9  = 19
    10 = 11
    11 = 18
    12 = 12
    13 = 14
    14 = 17
    15 = 8
    16 = 13
    17 = 9
    18 = 15
    19 = 16
    20 = 22
    21 = 23
    22 = 21
    23 = 20
}

$firstRow = 8
$lastCol = 51  # column AY

foreach ($oldRow in 8..23) {
    $newRow = $rowMap[$oldRow]
    $oldIdx = $oldRow - $firstRow + 1
    $newIdx = $newRow - $firstRow + 1

    for ($c = 1; $c -le $lastCol; $c++) {
        $v = $data[$oldIdx, $c]
        $ws.Cells.Item($newRow, $c).Value = $v
    }
}
